$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update machine id / name / mac / serial / ip for existing + new rows ---
$ws.Cells.Item(2, 1).Value = 10001
$ws.Cells.Item(2, 2).Value = "Machine 1"
$ws.Cells.Item(2, 3).Value = "8C-16-45-5A-5D-0D"
$ws.Cells.Item(2, 4).Value = "NM5328114630"
$ws.Cells.Item(2, 5).Value = "192.168.0.150"
$ws.Cells.Item(2, 6).Value = 1001
$ws.Cells.Item(2, 7).Value = "eng"
$ws.Cells.Item(2, 8).Value = $true
$ws.Cells.Item(2, 9).Value = "superadmin"
$ws.Cells.Item(2, 10).Value = "now()"

$ws.Cells.Item(3, 1).Value = 10002
$ws.Cells.Item(3, 2).Value = "Machine 2"
$ws.Cells.Item(3, 3).Value = "8C-16-45-88-E1-0D"
$ws.Cells.Item(3, 4).Value = "WY2132605316"
$ws.Cells.Item(3, 5).Value = "192.168.0.133"
$ws.Cells.Item(3, 6).Value = 1001
$ws.Cells.Item(3, 7).Value = "eng"
$ws.Cells.Item(3, 8).Value = $true
$ws.Cells.Item(3, 9).Value = "superadmin"
$ws.Cells.Item(3, 10).Value = "now()"

$ws.Cells.Item(4, 1).Value = 10003
$ws.Cells.Item(4, 2).Value = "Machine 3"
$ws.Cells.Item(4, 3).Value = "00-FF-D3-E3-9A-27"
$ws.Cells.Item(4, 4).Value = "CM6384145127"
$ws.Cells.Item(4, 5).Value = "192.168.0.161"
$ws.Cells.Item(4, 6).Value = 1001
$ws.Cells.Item(4, 7).Value = "eng"
$ws.Cells.Item(4, 8).Value = $true
$ws.Cells.Item(4, 9).Value = "superadmin"
$ws.Cells.Item(4, 10).Value = "now()"

$ws.Cells.Item(5, 1).Value = 10004
$ws.Cells.Item(5, 2).Value = "Machine 4"
$ws.Cells.Item(5, 3).Value = "8C-16-45-5A-62-41"
$ws.Cells.Item(5, 4).Value = "NT894252578"
$ws.Cells.Item(5, 5).Value = "192.168.0.259"
$ws.Cells.Item(5, 6).Value = 1001
$ws.Cells.Item(5, 7).Value = "eng"
$ws.Cells.Item(5, 8).Value = $true
$ws.Cells.Item(5, 9).Value = "superadmin"
$ws.Cells.Item(5, 10).Value = "now()"

$ws.Cells.Item(6, 1).Value = 10005
$ws.Cells.Item(6, 2).Value = "Machine 5"
$ws.Cells.Item(6, 3).Value = "E8-6A-64-1D-75-E4"
$ws.Cells.Item(6, 4).Value = "YM866672706"
$ws.Cells.Item(6, 5).Value = "192.168.0.119"
$ws.Cells.Item(6, 6).Value = 1001
$ws.Cells.Item(6, 7).Value = "eng"
$ws.Cells.Item(6, 8).Value = $true
$ws.Cells.Item(6, 9).Value = "superadmin"
$ws.Cells.Item(6, 10).Value = "now()"

$ws.Cells.Item(7, 1).Value = 10006
$ws.Cells.Item(7, 2).Value = "Machine 6"
$ws.Cells.Item(7, 3).Value = "8C-16-45-FA-94-B7"
$ws.Cells.Item(7, 4).Value = "WT6501645780"
$ws.Cells.Item(7, 5).Value = "192.168.0.177"
$ws.Cells.Item(7, 6).Value = 1001
$ws.Cells.Item(7, 7).Value = "eng"
$ws.Cells.Item(7, 8).Value = $true
$ws.Cells.Item(7, 9).Value = "superadmin"
$ws.Cells.Item(7, 10).Value = "now()"

$ws.Cells.Item(8, 1).Value = 10007
$ws.Cells.Item(8, 2).Value = "Machine 7"
$ws.Cells.Item(8, 3).Value = "8C-16-45-1A-0F-62"
$ws.Cells.Item(8, 4).Value = "LK8186452621"
$ws.Cells.Item(8, 5).Value = "192.168.0.227"
$ws.Cells.Item(8, 6).Value = 1001
$ws.Cells.Item(8, 7).Value = "eng"
$ws.Cells.Item(8, 8).Value = $true
$ws.Cells.Item(8, 9).Value = "superadmin"
$ws.Cells.Item(8, 10).Value = "now()"

$ws.Cells.Item(9, 1).Value = 10008
$ws.Cells.Item(9, 2).Value = "Machine 8"
$ws.Cells.Item(9, 3).Value = "E8-6A-64-1C-52-6E"
$ws.Cells.Item(9, 4).Value = "NR3264783870"
$ws.Cells.Item(9, 5).Value = "192.168.0.207"
$ws.Cells.Item(9, 6).Value = 1001
$ws.Cells.Item(9, 7).Value = "eng"
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(9, 9).Value = "superadmin"
$ws.Cells.Item(9, 10).Value = "now()"

$ws.Cells.Item(10, 1).Value = 10009
$ws.Cells.Item(10, 2).Value = "Machine 9"
$ws.Cells.Item(10, 3).Value = "48-51-B7-10-35-A6"
$ws.Cells.Item(10, 4).Value = "RW437027336"
$ws.Cells.Item(10, 5).Value = "192.168.0.220"
$ws.Cells.Item(10, 6).Value = 1001
$ws.Cells.Item(10, 7).Value = "eng"
$ws.Cells.Item(10, 8).Value = $true
$ws.Cells.Item(10, 9).Value = "superadmin"
$ws.Cells.Item(10, 10).Value = "now()"

$ws.Cells.Item(11, 1).Value = 10010
$ws.Cells.Item(11, 2).Value = "Machine 10"
$ws.Cells.Item(11, 3).Value = "8C-16-45-38-F3-F3"
$ws.Cells.Item(11, 4).Value = "SI158158531"
$ws.Cells.Item(11, 5).Value = "192.168.0.242"
$ws.Cells.Item(11, 6).Value = 1001
$ws.Cells.Item(11, 7).Value = "eng"
$ws.Cells.Item(11, 8).Value = $true
$ws.Cells.Item(11, 9).Value = "superadmin"
$ws.Cells.Item(11, 10).Value = "now()"

$ws.Cells.Item(12, 1).Value = 10011
$ws.Cells.Item(12, 2).Value = "Machine 11"
$ws.Cells.Item(12, 3).Value = "D4-3D-7E-58-CC-45"
$ws.Cells.Item(12, 4).Value = "XF3416823469"
$ws.Cells.Item(12, 5).Value = "192.168.0.173"
$ws.Cells.Item(12, 6).Value = 1001
$ws.Cells.Item(12, 7).Value = "eng"
$ws.Cells.Item(12, 8).Value = $true
$ws.Cells.Item(12, 9).Value = "superadmin"
$ws.Cells.Item(12, 10).Value = "now()"

$ws.Cells.Item(13, 1).Value = 10012
$ws.Cells.Item(13, 2).Value = "Machine 12"
$ws.Cells.Item(13, 3).Value = "8C-16-45-5A-5D-96"
$ws.Cells.Item(13, 4).Value = "BW4524978011"
$ws.Cells.Item(13, 5).Value = "192.168.0.203"
$ws.Cells.Item(13, 6).Value = 1001
$ws.Cells.Item(13, 7).Value = "eng"
$ws.Cells.Item(13, 8).Value = $true
$ws.Cells.Item(13, 9).Value = "superadmin"
$ws.Cells.Item(13, 10).Value = "now()"

$ws.Cells.Item(14, 1).Value = 10013
$ws.Cells.Item(14, 2).Value = "Machine 13"
$ws.Cells.Item(14, 3).Value = "8C-16-45-5A-5D-8E"
$ws.Cells.Item(14, 4).Value = "DB289579153"
$ws.Cells.Item(14, 5).Value = "192.168.0.112"
$ws.Cells.Item(14, 6).Value = 1001
$ws.Cells.Item(14, 7).Value = "eng"
$ws.Cells.Item(14, 8).Value = $true
$ws.Cells.Item(14, 9).Value = "superadmin"
$ws.Cells.Item(14, 10).Value = "now()"

$ws.Cells.Item(15, 1).Value = 10014
$ws.Cells.Item(15, 2).Value = "Machine 14"
$ws.Cells.Item(15, 3).Value = "8C-16-45-33-A5-5F"
$ws.Cells.Item(15, 4).Value = "SI4597903231"
$ws.Cells.Item(15, 5).Value = "192.168.0.178"
$ws.Cells.Item(15, 6).Value = 1001
$ws.Cells.Item(15, 7).Value = "eng"
$ws.Cells.Item(15, 8).Value = $true
$ws.Cells.Item(15, 9).Value = "superadmin"
$ws.Cells.Item(15, 10).Value = "now()"

$ws.Cells.Item(16, 1).Value = 10015
$ws.Cells.Item(16, 2).Value = "Machine 15"
$ws.Cells.Item(16, 3).Value = "3C-95-09-F9-EA-DF"
$ws.Cells.Item(16, 4).Value = "TJ7809002958"
$ws.Cells.Item(16, 5).Value = "192.168.0.267"
$ws.Cells.Item(16, 6).Value = 1001
$ws.Cells.Item(16, 7).Value = "eng"
$ws.Cells.Item(16, 8).Value = $true
$ws.Cells.Item(16, 9).Value = "superadmin"
$ws.Cells.Item(16, 10).Value = "now()"

$ws.Cells.Item(17, 1).Value = 10016
$ws.Cells.Item(17, 2).Value = "Machine 16"
$ws.Cells.Item(17, 3).Value = "8C-16-45-88-E7-0B"
$ws.Cells.Item(17, 4).Value = "JR6082789079"
$ws.Cells.Item(17, 5).Value = "192.168.0.149"
$ws.Cells.Item(17, 6).Value = 1001
$ws.Cells.Item(17, 7).Value = "eng"
$ws.Cells.Item(17, 8).Value = $true
$ws.Cells.Item(17, 9).Value = "superadmin"
$ws.Cells.Item(17, 10).Value = "now()"

$ws.Cells.Item(18, 1).Value = 10017
$ws.Cells.Item(18, 2).Value = "Machine 17"
$ws.Cells.Item(18, 3).Value = "B4-69-21-5A-DB-C4"
$ws.Cells.Item(18, 4).Value = "SA3722889241"
$ws.Cells.Item(18, 5).Value = "192.168.0.127"
$ws.Cells.Item(18, 6).Value = 1001
$ws.Cells.Item(18, 7).Value = "eng"
$ws.Cells.Item(18, 8).Value = $true
$ws.Cells.Item(18, 9).Value = "superadmin"
$ws.Cells.Item(18, 10).Value = "now()"

$ws.Cells.Item(19, 1).Value = 10018
$ws.Cells.Item(19, 2).Value = "Machine 18"
$ws.Cells.Item(19, 3).Value = "E8-6A-64-1D-48-B7"
$ws.Cells.Item(19, 4).Value = "RR2683722548"
$ws.Cells.Item(19, 5).Value = "192.168.0.248"
$ws.Cells.Item(19, 6).Value = 1001
$ws.Cells.Item(19, 7).Value = "eng"
$ws.Cells.Item(19, 8).Value = $true
$ws.Cells.Item(19, 9).Value = "superadmin"
$ws.Cells.Item(19, 10).Value = "now()"

$ws.Cells.Item(20, 1).Value = 10019
$ws.Cells.Item(20, 2).Value = "Machine 19"
$ws.Cells.Item(20, 3).Value = "8C-16-45-59-69-09 "
$ws.Cells.Item(20, 4).Value = "PO6528391346"
$ws.Cells.Item(20, 5).Value = "192.168.0.121"
$ws.Cells.Item(20, 6).Value = 1001
$ws.Cells.Item(20, 7).Value = "eng"
$ws.Cells.Item(20, 8).Value = $true
$ws.Cells.Item(20, 9).Value = "superadmin"
$ws.Cells.Item(20, 10).Value = "now()"

$ws.Cells.Item(21, 1).Value = 10020
$ws.Cells.Item(21, 2).Value = "Machine 20"
$ws.Cells.Item(21, 3).Value = "98-E7-F4-30-16-5A "
$ws.Cells.Item(21, 4).Value = "FB5962911652"
$ws.Cells.Item(21, 5).Value = "192.168.0.215"
$ws.Cells.Item(21, 6).Value = 1001
$ws.Cells.Item(21, 7).Value = "eng"
$ws.Cells.Item(21, 8).Value = $true
$ws.Cells.Item(21, 9).Value = "superadmin"
$ws.Cells.Item(21, 10).Value = "now()"

$ws.Cells.Item(22, 1).Value = 10021
$ws.Cells.Item(22, 2).Value = "Machine 21"
$ws.Cells.Item(22, 3).Value = "38-BA-F8-53-C7-8F"
$ws.Cells.Item(22, 4).Value = "FB5962911653"
$ws.Cells.Item(22, 5).Value = "192.168.0.874"
$ws.Cells.Item(22, 6).Value = 1001
$ws.Cells.Item(22, 7).Value = "eng"
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(22, 9).Value = "superadmin"
$ws.Cells.Item(22, 10).Value = "now()"

$ws.Cells.Item(23, 1).Value = 10022
$ws.Cells.Item(23, 2).Value = "Machine 22"
$ws.Cells.Item(23, 3).Value = "E8-6A-64-1C-58-C2"
$ws.Cells.Item(23, 4).Value = "FB5962911654"
$ws.Cells.Item(23, 5).Value = "192.168.0.721"
$ws.Cells.Item(23, 6).Value = 1001
$ws.Cells.Item(23, 7).Value = "eng"
$ws.Cells.Item(23, 8).Value = $true
$ws.Cells.Item(23, 9).Value = "superadmin"
$ws.Cells.Item(23, 10).Value = "now()"

$ws.Cells.Item(24, 1).Value = 10023
$ws.Cells.Item(24, 2).Value = "Machine 23"
$ws.Cells.Item(24, 3).Value = "E4-A4-71-CE-BA-93"
$ws.Cells.Item(24, 4).Value = "FB5962911655"
$ws.Cells.Item(24, 5).Value = "192.168.0.841"
$ws.Cells.Item(24, 6).Value = 1001
$ws.Cells.Item(24, 7).Value = "eng"
$ws.Cells.Item(24, 8).Value = $true
$ws.Cells.Item(24, 9).Value = "superadmin"
$ws.Cells.Item(24, 10).Value = "now()"

$ws.Cells.Item(25, 1).Value = 10024
$ws.Cells.Item(25, 2).Value = "Machine 24"
$ws.Cells.Item(25, 3).Value = "54-E1-AD-EA-30-C9"
$ws.Cells.Item(25, 4).Value = "FB5962911656"
$ws.Cells.Item(25, 5).Value = "192.168.0.186"
$ws.Cells.Item(25, 6).Value = 1001
$ws.Cells.Item(25, 7).Value = "eng"
$ws.Cells.Item(25, 8).Value = $true
$ws.Cells.Item(25, 9).Value = "superadmin"
$ws.Cells.Item(25, 10).Value = "now()"

$ws.Cells.Item(26, 1).Value = 10025
$ws.Cells.Item(26, 2).Value = "Machine 25"
$ws.Cells.Item(26, 3).Value = "8C-16-45-65-DD-40"
$ws.Cells.Item(26, 4).Value = "FB5962911657"
$ws.Cells.Item(26, 5).Value = "192.168.0.627"
$ws.Cells.Item(26, 6).Value = 1001
$ws.Cells.Item(26, 7).Value = "eng"
$ws.Cells.Item(26, 8).Value = $true
$ws.Cells.Item(26, 9).Value = "superadmin"
$ws.Cells.Item(26, 10).Value = "now()"

$ws.Cells.Item(27, 1).Value = 10026
$ws.Cells.Item(27, 2).Value = "Machine 26"
$ws.Cells.Item(27, 3).Value = "58-20-B1-D6-C3-BE"
$ws.Cells.Item(27, 4).Value = "FB5962911658"
$ws.Cells.Item(27, 5).Value = "192.168.0.879"
$ws.Cells.Item(27, 6).Value = 1001
$ws.Cells.Item(27, 7).Value = "eng"
$ws.Cells.Item(27, 8).Value = $true
$ws.Cells.Item(27, 9).Value = "superadmin"
$ws.Cells.Item(27, 10).Value = "now()"

$ws.Cells.Item(28, 1).Value = 10027
$ws.Cells.Item(28, 2).Value = "Machine 27"
$ws.Cells.Item(28, 3).Value = "8C-16-45-38-F0-25"
$ws.Cells.Item(28, 4).Value = "FB5962911659"
$ws.Cells.Item(28, 5).Value = "192.168.0.628"
$ws.Cells.Item(28, 6).Value = 1001
$ws.Cells.Item(28, 7).Value = "eng"
$ws.Cells.Item(28, 8).Value = $true
$ws.Cells.Item(28, 9).Value = "superadmin"
$ws.Cells.Item(28, 10).Value = "now()"

$ws.Cells.Item(29, 1).Value = 10028
$ws.Cells.Item(29, 2).Value = "Machine 28"
$ws.Cells.Item(29, 3).Value = "6C-88-14-AC-EF-55"
$ws.Cells.Item(29, 4).Value = "FB5962911661"
$ws.Cells.Item(29, 5).Value = "192.168.0.306"
$ws.Cells.Item(29, 6).Value = 1001
$ws.Cells.Item(29, 7).Value = "eng"
$ws.Cells.Item(29, 8).Value = $true
$ws.Cells.Item(29, 9).Value = "superadmin"
$ws.Cells.Item(29, 10).Value = "now()"

$ws.Cells.Item(30, 1).Value = 10029
$ws.Cells.Item(30, 2).Value = "Machine 29"
$ws.Cells.Item(30, 3).Value = "3C-6A-A7-C0-DF-27"
$ws.Cells.Item(30, 4).Value = "FB5962911662"
$ws.Cells.Item(30, 5).Value = "192.168.0.355"
$ws.Cells.Item(30, 6).Value = 1001
$ws.Cells.Item(30, 7).Value = "eng"
$ws.Cells.Item(30, 8).Value = $true
$ws.Cells.Item(30, 9).Value = "superadmin"
$ws.Cells.Item(30, 10).Value = "now()"

# --- Widen mac_address column to fit the new, longer hyphenated MAC format ---
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668

# --- Move selection to the first empty row below the table (whole-row selection) ---
$ws.Range("A31:XFD1048576").Select() | Out-Null
